$d = $word.ActiveDocument

$replacements = @(
    @("2024-09-23 Monday", "2024-09-24 Tuesday"),
    @("423÷8=52, 7", "638÷3=212, 2"),
    @("430÷5=86, 0", "227÷5=45, 2"),
    @("526÷6=87, 4", "659÷5=131, 4"),
    @("809÷6=134, 5", "387÷6=64, 3"),
    @("162÷6=27, 0", "957÷6=159, 3"),
    @("269÷6=44, 5", "256÷2=128, 0"),
    @("728÷4=182, 0", "746÷9=82, 8"),
    @("626÷4=156, 2", "801÷2=400, 1"),
    @("944÷4=236, 0", "120÷2=60, 0"),
    @("817÷2=408, 1", "586÷7=83, 5"),
    @("925÷6=154, 1", "898÷3=299, 1"),
    @("196÷4=49, 0", "730÷7=104, 2"),
    @("381÷8=47, 5", "801÷7=114, 3"),
    @("528÷3=176, 0", "670÷6=111, 4"),
    @("336÷8=42, 0", "724÷7=103, 3"),
    @("429÷6=71, 3", "948÷9=105, 3"),
    @("941÷7=134, 3", "502÷2=251, 0"),
    @("383÷7=54, 5", "184÷6=30, 4"),
    @("384÷4=96, 0", "676÷6=112, 4"),
    @("199÷5=39, 4", "390÷3=130, 0"),
    @("763÷4=190, 3", "733÷7=104, 5"),
    @("218÷4=54, 2", "771÷5=154, 1"),
    @("165÷6=27, 3", "217÷4=54, 1"),
    @("367÷2=183, 1", "733÷6=122, 1"),
    @("316÷7=45, 1", "108÷6=18, 0")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
